# Weekly update: insert 3 new price rows (Clemenuless, "Terminal La Palmera
# de La Serena" / Mandarina sheet) above the existing data, pushing the
# previously-existing rows 560:584 down to 563:587.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 560, shifting the old 560:584 block down to 563:587.
$ws.Rows("560:562").Insert()

# New row 560: Clemenuless / Especial
$ws.Cells.Item(560, 1).Value = 8
$ws.Cells.Item(560, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(560, 3).Value = "Coquimbo"
$ws.Cells.Item(560, 4).Value = 44753
$ws.Cells.Item(560, 5).Value = 4
$ws.Cells.Item(560, 6).Value = "Fruta"
$ws.Cells.Item(560, 7).Value = 100102
$ws.Cells.Item(560, 8).Value = "Cítricos"
$ws.Cells.Item(560, 9).Value = 100102004
$ws.Cells.Item(560, 10).Value = "Mandarina"
$ws.Cells.Item(560, 11).Value = "Clemenuless"
$ws.Cells.Item(560, 12).Value = "Especial"
$ws.Cells.Item(560, 13).Value = 400
$ws.Cells.Item(560, 14).Value = 6500
$ws.Cells.Item(560, 15).Value = 7000
$ws.Cells.Item(560, 16).Value = 6750
$ws.Cells.Item(560, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(560, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(560, 19).Value = 675
$ws.Cells.Item(560, 20).Value = 10

# New row 561: Clemenuless / Primera
$ws.Cells.Item(561, 1).Value = 8
$ws.Cells.Item(561, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(561, 3).Value = "Coquimbo"
$ws.Cells.Item(561, 4).Value = 44753
$ws.Cells.Item(561, 5).Value = 4
$ws.Cells.Item(561, 6).Value = "Fruta"
$ws.Cells.Item(561, 7).Value = 100102
$ws.Cells.Item(561, 8).Value = "Cítricos"
$ws.Cells.Item(561, 9).Value = 100102004
$ws.Cells.Item(561, 10).Value = "Mandarina"
$ws.Cells.Item(561, 11).Value = "Clemenuless"
$ws.Cells.Item(561, 12).Value = "Primera"
$ws.Cells.Item(561, 13).Value = 500
$ws.Cells.Item(561, 14).Value = 4500
$ws.Cells.Item(561, 15).Value = 5000
$ws.Cells.Item(561, 16).Value = 4750
$ws.Cells.Item(561, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(561, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(561, 19).Value = 475
$ws.Cells.Item(561, 20).Value = 10

# New row 562: Clemenuless / Segunda
$ws.Cells.Item(562, 1).Value = 8
$ws.Cells.Item(562, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(562, 3).Value = "Coquimbo"
$ws.Cells.Item(562, 4).Value = 44753
$ws.Cells.Item(562, 5).Value = 4
$ws.Cells.Item(562, 6).Value = "Fruta"
$ws.Cells.Item(562, 7).Value = 100102
$ws.Cells.Item(562, 8).Value = "Cítricos"
$ws.Cells.Item(562, 9).Value = 100102004
$ws.Cells.Item(562, 10).Value = "Mandarina"
$ws.Cells.Item(562, 11).Value = "Clemenuless"
$ws.Cells.Item(562, 12).Value = "Segunda"
$ws.Cells.Item(562, 13).Value = 360
$ws.Cells.Item(562, 14).Value = 2500
$ws.Cells.Item(562, 15).Value = 3000
$ws.Cells.Item(562, 16).Value = 2750
$ws.Cells.Item(562, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(562, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(562, 19).Value = 275
$ws.Cells.Item(562, 20).Value = 10
